$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

# Column K = "Stock Mínimo Objetivo", Column L = "Diferencia Stock"
# Set L(row) = K(row) for data rows 3..132
for ($r = 3; $r -le 132; $r++) {
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 12).Value2 = $kVal
}

# Update the summary total "Total_Ajuste_Stock" (C146) to the sum of the
# Diferencia Stock column (L3:L132), mirroring Stock_Minimo_Objetivo (C145)
$total = 0
for ($r = 3; $r -le 132; $r++) {
    $total += $ws.Cells.Item($r, 12).Value2
}
$ws.Range("C146").Value2 = $total
